# Avance para sexta entrega
# Update the F-column numeric values on the "Casos de Uso" sheet and
# move the sheet's visible/selected cell so the view reflects the last
# edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Update F-column values (+3 offsets applied by the author's edit)
$ws.Range("F8").Value = 3.37
$ws.Range("F24").Value = 5.58
$ws.Range("F26").Value = 6.16
$ws.Range("F28").Value = 5.76
$ws.Range("F30").Value = 15.63
$ws.Range("F32").Value = 9.53
$ws.Range("F36").Value = 10

# Update the sheet view: scroll so D24 is the top-left visible cell and
# select F32, matching the saved view state in the workbook.
$ws.Activate()
$ws.Range("F32").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 4
